$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ybus formation script: recompute the admittance matrix.
# Narrow the Bus4/Bus5/Bus6 columns (E, F, G) to line up with the other
# bus columns (H, I, J were already narrower).
$ws.Columns("E").ColumnWidth = 5.333333
$ws.Columns("F").ColumnWidth = 5.333333
$ws.Columns("G").ColumnWidth = 5.333333

# Updated Ybus admittance values produced by the new formation script.
$ws.Range("E5").Value = 3.307
$ws.Range("F5").Value = -1.365
$ws.Range("G5").Value = -1.942
$ws.Range("F6").Value = 2.553
$ws.Range("I9").Value = 2.772
